# docs/datadictionary.xlsx — "added data screenshots, checked data dictionary for typos"
#
# Data-dictionary grid on Sheet1 had two typos in its "Field name" column:
#   B1  "Filed Name" -> "Field name"
#   A23 "Tranc2"     -> "Trans2"
# Fix them in place (this reorders the shared-strings table exactly like the
# authored diff: the two old strings drop out and the corrected ones are
# appended at the end of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Field name"
$ws.Range("A23").Value = "Trans2"

# Leave the selection where the author left it after reviewing the fixes.
$ws.Range("A24").Select() | Out-Null
